# Updated symbol list on Sat Feb 11 17:25:56 UTC 2023 with GitHub Actions
# Refreshes the crypto price/volume snapshot (columns D/E) and shifts the
# Coin/Link rows 6-17 down one slot to make room for the new "GateToken"
# entry. Numeric-looking text values are written with a leading apostrophe
# so Excel stores them as text (matching the sheet's existing inlineStr
# "number as text" cells) instead of converting them to real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.64"
$ws.Range("E2").Value = "'0.74%"
$ws.Range("D3").Value = "'40.80"
$ws.Range("E3").Value = "'1.00%"
$ws.Range("D4").Value = "'5.127"
$ws.Range("E4").Value = "'1.51%"
$ws.Range("D5").Value = "'0.07621"
$ws.Range("E5").Value = "'0.40%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.279"
$ws.Range("E6").Value = "'0.73%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.623"
$ws.Range("E7").Value = "'2.22%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.499"
$ws.Range("E8").Value = "'2.24%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9093"
$ws.Range("E9").Value = "'0.49%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1180"
$ws.Range("E10").Value = "'17.94%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1817"
$ws.Range("E11").Value = "'3.74%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09175"
$ws.Range("E12").Value = "'2.03%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04260"
$ws.Range("E13").Value = "'-3.02%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1045"
$ws.Range("E14").Value = "'-0.80%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001250"
$ws.Range("E15").Value = "'-0.54%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005797"
$ws.Range("E16").Value = "'-1.13%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.358"
$ws.Range("E17").Value = "'-0.35%"
$ws.Range("E19").Value = "'0.88%"
$ws.Range("E20").Value = "'4.43%"
$ws.Range("D21").Value = "'0.2706"
$ws.Range("E21").Value = "'-5.09%"
$ws.Range("D22").Value = "'0.04048"
$ws.Range("E22").Value = "'-2.69%"
$ws.Range("E23").Value = "'4.33%"
$ws.Range("D24").Value = "'0.004108"
$ws.Range("E24").Value = "'1.23%"
$ws.Range("D25").Value = "'0.0001272"
$ws.Range("E25").Value = "'-2.48%"
$ws.Range("D26").Value = "'0.0003751"
$ws.Range("D38").Value = "'0.02424"
$ws.Range("E38").Value = "'1.13%"
$ws.Range("D39").Value = "'0.05228"
$ws.Range("E39").Value = "'2.05%"
$ws.Range("D40").Value = "'0.007782"
$ws.Range("E40").Value = "'-0.93%"
$ws.Range("D41").Value = "'0.1299"
$ws.Range("E41").Value = "'-0.12%"
$ws.Range("D42").Value = "'0.006804"
$ws.Range("E42").Value = "'-4.35%"
$ws.Range("D43").Value = "'0.001933"
$ws.Range("E43").Value = "'-1.11%"
$ws.Range("D44").Value = "'0.008087"
$ws.Range("E44").Value = "'-3.47%"
$ws.Range("E45").Value = "'-7.62%"
$ws.Range("D46").Value = "'0.00006898"
$ws.Range("E46").Value = "'6.83%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.15%"
$ws.Range("D48").Value = "'0.09626"
$ws.Range("E48").Value = "'1,668.48%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.15%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'-0.15%"
